$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "CFP10 SR10" was renamed to "CFP2 SR10" (part-number label column A)
$ws.Range("A9").Value = "CFP2 SR10"

# New rows appended at the bottom of the table (34-40). Labels for the first
# two new rows were typed down column A first ...
$ws.Range("A34").Value = "FDR"
$ws.Range("A35").Value = "FDR Gen3"

# ... then the matching part numbers filled into column C ...
$ws.Range("C34").Value = "FCBN414QB1;FCBG414QB1"
$ws.Range("C35").Value = "FCBN414QD3;FCCN414QD3"

# Row 19 (Quadwire Gen3): the part-number list in column C was split out;
# this sheet keeps only the shorter/no-QD414 variant now.
$ws.Range("C19").Value = "FCCG410QD3;FCBG410QD3;FCBN410QD3;FCCN410QD3"

# Remaining new rows filled in A/C pairs, one row at a time.
$ws.Range("A36").Value = "FDR transiver"
$ws.Range("C36").Value = "FTL414QB2;FTL414QL2"

$ws.Range("A37").Value = "SNAP12"
$ws.Range("C37").Value = "FTXD02SL1"

$ws.Range("A38").Value = "Octopus"
$ws.Range("C38").Value = "FCBR510QE2;FCBN510QE2"

$ws.Range("A39").Value = "QSFPSR4 FET Gen2"
$ws.Range("C39").Value = "FTL410QT2"

$ws.Range("A40").Value = "QSFPSR4 FET Gen3"
$ws.Range("C40").Value = "FTL410QT3"

# Update the view: active cell/selection moved to the newly added last row,
# and the window scrolled down so row 19 is at the top.
$ws.Range("C40").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
